$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.510.26"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "1.681.51"
$ws.Range("E3").Value = "  +2.47%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.48"
$ws.Range("E5").Value = "  +4.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5316"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2673"
$ws.Range("E8").Value = "  +4.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06413"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.44"
$ws.Range("E10").Value = "  +5.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").Value = "1.689.75"
$ws.Range("E12").Value = "  +3.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.508"
$ws.Range("E13").Value = "  +3.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5617"
$ws.Range("E14").Value = "  +4.37%  "
$ws.Range("D15").Value = "0.0₅8392"
$ws.Range("E15").Value = "  +5.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.91"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("D17").Value = "26.544.49"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.809"
$ws.Range("E19").Value = "  +3.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.65"
$ws.Range("E20").Value = "  +4.90%  "
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.392"
$ws.Range("E22").Value = "  +5.03%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.41"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1272"
$ws.Range("E25").Value = "  +6.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.471"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  +5.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.413"
$ws.Range("E28").Value = "  +3.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06127"
$ws.Range("E29").Value = "  +2.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.278"
$ws.Range("E30").Value = "  +2.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.611"
$ws.Range("E31").Value = "  +7.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.459"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.702"
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.014"
$ws.Range("E34").Value = "  +4.83%  "
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5705"
$ws.Range("E37").Value = "  -2.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01640"
$ws.Range("E38").Value = "  +3.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.962"
$ws.Range("E39").Value = "  +3.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8668"
$ws.Range("E40").Value = "  +3.27%  "
$ws.Range("D41").Value = "1.058.49"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.08"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").Value = "1.831.64"
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.20"
$ws.Range("E46").Value = "  +5.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.128"
$ws.Range("E47").Value = "  +2.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9997"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05206"
$ws.Range("E49").Value = "  +0.23%  "
$ws.Range("E50").Value = "  +4.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4243"
$ws.Range("E51").Value = "  +0.30%  "
